$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H43").Value = 1528.421
$ws.Range("I43").Value = 2235.375
$ws.Range("J43").Value = 1014.2727
$ws.Range("K43").Value = 2235.375
$ws.Range("L43").Value = 1014.2727
$ws.Range("M43").Value = -2166.375
$ws.Range("N43").Value = -1152.2727
$ws.Range("H51").Value = 9820.066000000001
$ws.Range("I51").Value = 27475
$ws.Range("J51").Value = 3400.0908
$ws.Range("K51").Value = 27475
$ws.Range("L51").Value = 3400.0908
$ws.Range("M51").Value = -26991
$ws.Range("N51").Value = -4368.0908
$ws.Range("H64").Value = 29189.078
$ws.Range("I64").Value = 65181.188
$ws.Range("J64").Value = 3013
$ws.Range("K64").Value = 65181.188
$ws.Range("L64").Value = 3013
$ws.Range("M64").Value = -64933.188
$ws.Range("N64").Value = -3509
$ws.Range("H67").Value = 29189.078
$ws.Range("I67").Value = 65181.188
$ws.Range("J67").Value = 3013
$ws.Range("K67").Value = 65181.188
$ws.Range("L67").Value = 3013
$ws.Range("M67").Value = -64323.188
$ws.Range("N67").Value = -4729
$ws.Range("H116").Value = 3161.2
$ws.Range("I116").Value = 1980
$ws.Range("K116").Value = 1980
$ws.Range("M116").Value = 1462
$ws.Range("H127").Value = 58825036
$ws.Range("I127").Value = 549.25
$ws.Range("J127").Value = 76924880
$ws.Range("K127").Value = 1647.75
$ws.Range("L127").Value = 230774640
$ws.Range("M127").Value = 3312.25
$ws.Range("N127").Value = -230784560
$ws.Range("H129").Value = 2308.3809
$ws.Range("J129").Value = 799.625
$ws.Range("L129").Value = 2398.875
$ws.Range("N129").Value = -12398.875
$ws.Range("H137").Value = 1254.9354
$ws.Range("I137").Value = 1340.3636
$ws.Range("J137").Value = 1046.1111
$ws.Range("K137").Value = 4021.0908
$ws.Range("L137").Value = 3138.3333
$ws.Range("M137").Value = -1471.0908
$ws.Range("N137").Value = -8238.3333

$ws = $wb.Worksheets("ARM")
$ws.Range("H61").Value = 2243.4595
$ws.Range("I61").Value = 2023.2963
$ws.Range("K61").Value = 2023.2963
$ws.Range("M61").Value = -1811.2963
$ws.Range("H97").Value = 28778.277
$ws.Range("I97").Value = 40700.36
$ws.Range("J97").Value = 1682.6364
$ws.Range("K97").Value = 40700.36
$ws.Range("L97").Value = 1682.6364
$ws.Range("M97").Value = -40204.36
$ws.Range("N97").Value = -2674.6364
$ws.Range("H118").Value = 39300
$ws.Range("J118").Value = 39300
$ws.Range("L118").Value = 39300
$ws.Range("N118").Value = -42614
$ws.Range("H122").Value = 2561.4
$ws.Range("I122").Value = 2437.6667
$ws.Range("K122").Value = 7313.000100000001
$ws.Range("M122").Value = -4863.000100000001
$ws.Range("H136").Value = 2243.4595
$ws.Range("I136").Value = 2023.2963
$ws.Range("K136").Value = 6069.8889
$ws.Range("M136").Value = -3519.8889

$ws = $wb.Worksheets("BSM")
$ws.Range("H20").Value = 42120.28
$ws.Range("I20").Value = 57939.555
$ws.Range("K20").Value = 57939.555
$ws.Range("M20").Value = -57692.555

$ws = $wb.Worksheets("CRP")
$ws.Range("H134").Value = 1960.1
$ws.Range("I134").Value = 1733.4445
$ws.Range("K134").Value = 5200.333500000001
$ws.Range("M134").Value = -2665.333500000001

$ws = $wb.Worksheets("CUL")
$ws.Range("H5").Value = 1272.25
$ws.Range("I5").Value = 1464.3
$ws.Range("J5").Value = 1184.9546
$ws.Range("K5").Value = 4392.9
$ws.Range("L5").Value = 3554.8638
$ws.Range("M5").Value = -4280.9
$ws.Range("N5").Value = -3778.8638
$ws.Range("H18").Value = 566
$ws.Range("I18").Value = 499.81818
$ws.Range("K18").Value = 1499.45454
$ws.Range("M18").Value = -1330.45454
$ws.Range("H86").Value = 840.4167
$ws.Range("I86").Value = 482
$ws.Range("J86").Value = 873
$ws.Range("K86").Value = 1446
$ws.Range("L86").Value = 2619
$ws.Range("M86").Value = -260
$ws.Range("N86").Value = -4991
$ws.Range("H89").Value = 840.4167
$ws.Range("I89").Value = 482
$ws.Range("J89").Value = 873
$ws.Range("K89").Value = 4338
$ws.Range("L89").Value = 7857
$ws.Range("M89").Value = 1590
$ws.Range("N89").Value = -19713
$ws.Range("H131").Value = 735.7
$ws.Range("I131").Value = 340
$ws.Range("J131").Value = 756.5263
$ws.Range("K131").Value = 1020
$ws.Range("L131").Value = 2269.5789
$ws.Range("M131").Value = 4020
$ws.Range("N131").Value = -12349.5789
$ws.Range("H135").Value = 1272.25
$ws.Range("I135").Value = 1464.3
$ws.Range("J135").Value = 1184.9546
$ws.Range("K135").Value = 13178.7
$ws.Range("L135").Value = 10664.5914
$ws.Range("M135").Value = -10643.7
$ws.Range("N135").Value = -15734.5914

$ws = $wb.Worksheets("GSM")
$ws.Range("H126").Value = 3924062.2
$ws.Range("I126").Value = 3000.375
$ws.Range("K126").Value = 9001.125
$ws.Range("M126").Value = -6531.125

$ws = $wb.Worksheets("LTW")
$ws.Range("H7").Value = 3433.3057
$ws.Range("J7").Value = 3681.4285
$ws.Range("L7").Value = 3681.4285
$ws.Range("N7").Value = -3905.4285
$ws.Range("H16").Value = 4251182.5
$ws.Range("I16").Value = 7000806
$ws.Range("J16").Value = 715951.8
$ws.Range("K16").Value = 7000806
$ws.Range("L16").Value = 715951.8
$ws.Range("M16").Value = -7000636
$ws.Range("N16").Value = -716291.8
$ws.Range("H22").Value = 2036
$ws.Range("I22").Value = 1773.1666
$ws.Range("J22").Value = 2193.7
$ws.Range("K22").Value = 1773.1666
$ws.Range("L22").Value = 2193.7
$ws.Range("M22").Value = -1478.1666
$ws.Range("N22").Value = -2783.7
$ws.Range("H27").Value = 2036
$ws.Range("I27").Value = 1773.1666
$ws.Range("J27").Value = 2193.7
$ws.Range("K27").Value = 1773.1666
$ws.Range("L27").Value = 2193.7
$ws.Range("M27").Value = -1666.1666
$ws.Range("N27").Value = -2407.7
$ws.Range("H40").Value = 68867
$ws.Range("I40").Value = 168417.5
$ws.Range("K40").Value = 168417.5
$ws.Range("M40").Value = -168281.5
$ws.Range("H61").Value = 1780.7778
$ws.Range("I61").Value = 1776.4
$ws.Range("K61").Value = 1776.4
$ws.Range("M61").Value = -1574.4
$ws.Range("H100").Value = 1972.3
$ws.Range("J100").Value = 2020.5
$ws.Range("L100").Value = 2020.5
$ws.Range("N100").Value = -3102.5
$ws.Range("H113").Value = 1780.7778
$ws.Range("I113").Value = 1776.4
$ws.Range("K113").Value = 1776.4
$ws.Range("M113").Value = 393.5999999999999
$ws.Range("H122").Value = 2990
$ws.Range("I122").Value = 2990
$ws.Range("K122").Value = 8970
$ws.Range("M122").Value = -6520
$ws.Range("H126").Value = 3433.3057
$ws.Range("J126").Value = 3681.4285
$ws.Range("L126").Value = 11044.2855
$ws.Range("N126").Value = -15984.2855
$ws.Range("H136").Value = 1766.2222
$ws.Range("I136").Value = 1760.125
$ws.Range("K136").Value = 5280.375
$ws.Range("M136").Value = -2730.375

$ws = $wb.Worksheets("WVR")
$ws.Range("H62").Value = 6995443
$ws.Range("I62").Value = 38462784
$ws.Range("J62").Value = 2700
$ws.Range("K62").Value = 38462784
$ws.Range("L62").Value = 2700
$ws.Range("M62").Value = -38462160
$ws.Range("N62").Value = -3948
$ws.Range("H65").Value = 6995443
$ws.Range("I65").Value = 38462784
$ws.Range("J65").Value = 2700
$ws.Range("K65").Value = 192313920
$ws.Range("L65").Value = 13500
$ws.Range("M65").Value = -192310800
$ws.Range("N65").Value = -19740
$ws.Range("H126").Value = 2157.7
$ws.Range("I126").Value = 2009.875
$ws.Range("J126").Value = 2749
$ws.Range("K126").Value = 6029.625
$ws.Range("L126").Value = 8247
$ws.Range("M126").Value = -3559.625
$ws.Range("N126").Value = -13187
$ws.Range("H133").Value = 45715
$ws.Range("J133").Value = 45715
$ws.Range("L133").Value = 45715
$ws.Range("N133").Value = -55835
$ws.Range("H136").Value = 1162.1207
$ws.Range("I136").Value = 432.75674
$ws.Range("J136").Value = 2447.1904
$ws.Range("K136").Value = 1298.27022
$ws.Range("L136").Value = 7341.5712
$ws.Range("M136").Value = 1251.72978
$ws.Range("N136").Value = -12441.5712
